# Insert two new rows at the top of the data block (539:540), shifting the
# existing rows (old 539..560) down to (541..562) - this matches the
# dimension change from A1:T560 to A1:T562.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("539:540").Insert()

# Row 539 - new "Especial" Murcott entry for 2021-11-16 (serial 44516)
$ws.Cells.Item(539, 1).Value  = 3
$ws.Cells.Item(539, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(539, 3).Value  = "Coquimbo"
$ws.Cells.Item(539, 4).Value  = "2021-11-16"
$ws.Cells.Item(539, 5).Value  = 5
$ws.Cells.Item(539, 6).Value  = "Fruta"
$ws.Cells.Item(539, 7).Value  = 100102
$ws.Cells.Item(539, 8).Value  = "Cítricos"
$ws.Cells.Item(539, 9).Value  = 100102004
$ws.Cells.Item(539, 10).Value = "Mandarina"
$ws.Cells.Item(539, 11).Value = "Murcott"
$ws.Cells.Item(539, 12).Value = "Especial"
$ws.Cells.Item(539, 13).Value = 65
$ws.Cells.Item(539, 14).Value = 6000
$ws.Cells.Item(539, 15).Value = 6000
$ws.Cells.Item(539, 16).Value = 6000
$ws.Cells.Item(539, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(539, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(539, 19).Value = 600
$ws.Cells.Item(539, 20).Value = 10

# Row 540 - new "Primera" Murcott entry for 2021-11-16 (serial 44516)
$ws.Cells.Item(540, 1).Value  = 3
$ws.Cells.Item(540, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(540, 3).Value  = "Coquimbo"
$ws.Cells.Item(540, 4).Value  = "2021-11-16"
$ws.Cells.Item(540, 5).Value  = 5
$ws.Cells.Item(540, 6).Value  = "Fruta"
$ws.Cells.Item(540, 7).Value  = 100102
$ws.Cells.Item(540, 8).Value  = "Cítricos"
$ws.Cells.Item(540, 9).Value  = 100102004
$ws.Cells.Item(540, 10).Value = "Mandarina"
$ws.Cells.Item(540, 11).Value = "Murcott"
$ws.Cells.Item(540, 12).Value = "Primera"
$ws.Cells.Item(540, 13).Value = 60
$ws.Cells.Item(540, 14).Value = 5000
$ws.Cells.Item(540, 15).Value = 5000
$ws.Cells.Item(540, 16).Value = 5000
$ws.Cells.Item(540, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(540, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(540, 19).Value = 500
$ws.Cells.Item(540, 20).Value = 10
